# Apply the "chap 3 3rd personed, chapman update" edit to the
# reviewer-response workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reviewer 1's comment header now credits the reviewer by name.
$ws.Range("A2").Value = "REVIEWER 1: M Hart"

# 2. Row 34 (candidate's response to the Chapter-1 chemistry-error comment
#    in A34) gets a RESPONSE and CHANGES TO THESIS entry explaining the
#    Chapman-cycle / photolysis wavelength correction.
$ws.Range("B34").Value = "Chemistry is not my strongest suit, so I do not pick up errors as readily as I should " + [char]0x2013 + " thanks for pointing these out. I had "
$ws.Range("C34").Value = "Low wavelengths discussed in equation set 1.2 are updated to 350nm, with reference updated"

# The extra wrapped text slightly changes the auto-fitted row height.
$ws.Rows.Item(34).RowHeight = 242.5

# 3. Move the active selection / viewport down to the row that was edited.
$ws.Range("B34").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
